$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values for rows 2-5 and row 8
$ws.Range("F2").Value = -4
$ws.Range("F3").Value = 7
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = -2
$ws.Range("F8").Value = -7
